# DRF_Tutorial_Steps.docx edit:
#   1. Insert a new "Video About Refresh Tokens" list item right before
#      the existing "Create Cars app" list item.
#   2. Move the <w:lastRenderedPageBreak/> marker from the run that holds
#      "Re-do functions to be view functions" to the run that holds the
#      preceding "Commit" (i.e. it now renders one item earlier).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: new bullet "Video About Refresh Tokens" before "Create Cars app"
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a") -eq "Create Cars app") {
        $target = $i
        break
    }
}

$createCarsPara = $d.Paragraphs.Item($target)
$null = $createCarsPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($target)
$newPara.Range.Text = "Video About Refresh Tokens"

# ---------------------------------------------------------------------
# Change 2: relocate <w:lastRenderedPageBreak/> from "Re-do functions to
# be view functions" onto the "Commit" paragraph immediately before it.
# ---------------------------------------------------------------------
$redoIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a") -eq "Re-do functions to be view functions") {
        $redoIdx = $i
        break
    }
}

$redoPara = $d.Paragraphs.Item($redoIdx)
$commitPara = $d.Paragraphs.Item($redoIdx - 1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Add the page-break marker to the "Commit" run (re-emit the paragraph,
# preserving its identity attributes, with the marker ahead of the text).
$null = $commitPara.Range.InsertXML('<w:p ' + $ns + ' w14:paraId="6ED6FEAA" w14:textId="4BE159C6" w:rsidR="000110AD" w:rsidRDefault="000110AD" w:rsidP="000A23AD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Commit</w:t></w:r></w:p>')

# Remove the page-break marker from "Re-do functions to be view functions".
$redoPara = $d.Paragraphs.Item($redoIdx)
$null = $redoPara.Range.InsertXML('<w:p ' + $ns + ' w14:paraId="1F3DD1C8" w14:textId="5D0EBCD5" w:rsidR="000110AD" w:rsidRDefault="000110AD" w:rsidP="000A23AD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Re-do functions to be view functions</w:t></w:r></w:p>')
